$wb = $excel.ActiveWorkbook

# README sheet: update the bill-recording instruction text (dropped the
# "senate=1 or house=0" numeric-coding note in favor of plain language).
$readme = $wb.Worksheets.Item("README")
$readme.Range("A2").Value = "Record the state, year, title, and house (senate or not) for each bill. Use the year that the bill was actually voted on, dates should usually be included. "

# Scores sheet: remove the unused "Party" and "Pictures" columns (G:H).
$scores = $wb.Worksheets.Item("Scores")
$scores.Columns("G:H").Delete()
